$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- View changes: zoom to 70%, and select D1 ---
$excel.ActiveWindow.Zoom = 70
$ws.Range("D1").Select()

# --- H15:H18 gain "Finanza Matematica" (copy format from the neighboring G column cell) ---
foreach ($r in 15..18) {
    $ws.Range("G$r").Copy()
    $ws.Range("H$r").PasteSpecial(-4122) # xlPasteFormats
    $ws.Range("H$r").Value = "Finanza Matematica"
}
$excel.CutCopyMode = $false

# --- Row 19 ---
$ws.Range("F19").Value = "Introduzione QuantLib"
$ws.Range("G19").Value = "1° esercitazione QuantLib Python"
$ws.Range("F19").Copy()
$ws.Range("H19").PasteSpecial(-4122) # xlPasteFormats
$ws.Range("H19").Value = "2° esercitazione QuantLib Python"
$excel.CutCopyMode = $false

# --- Row 20 ---
$ws.Range("D20").Value = "Git"
$ws.Range("F20").Value = "Introduzione QuantLib"
$ws.Range("G20").Value = "1° esercitazione QuantLib Python"
$ws.Range("F20").Copy()
$ws.Range("H20").PasteSpecial(-4122) # xlPasteFormats
$ws.Range("H20").Value = "2° esercitazione QuantLib Python"
$excel.CutCopyMode = $false

# --- E29 / E30 ---
$ws.Range("E29").Value = "Introduzione Python"
$ws.Range("E30").Value = "Introduzione Python"
